# "Updated Features with Add CI steps"
# Replace the stray "dk gdpr purpose test" row with a new CI entry
# (SE_AAA_TestMartinx), give it wrap-text formatting, and leave the
# selection on the newly touched rows, matching the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 49 (A49) becomes the new CI "SE_AAA_TestMartinx" with wrapped text.
$ws.Range("A49").Value = "SE_AAA_TestMartinx"
$ws.Range("A49").WrapText = $true

# Final selection left on B49:B50 (active cell B49), matching the saved view.
[void]$ws.Range("B49:B50").Select()
